# Generate Report for Handback
# Update the timestamp values recorded on the handback-status report sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date for the 23a6492f... entry,
# also shared with de-de!H2 (same literal timestamp text).
$wsOverview.Range("G2").Value = "2016-09-02 11:17:22"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-09-02 11:17:17"
$wsZhCn.Range("K2").Value = "2016-09-02 11:17:35"

# de-de sheet: Correspond Handoff Datetime (shared w/ Overview!G2) / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-02 11:17:22"
$wsDeDe.Range("K2").Value = "2016-09-02 11:17:42"
